{"js": "// The SDO template merge-field referencing the case management location's\n// venue name field should instead reference the location ref data's\n// \"external short name\" field: `venue_name` -> `external_short_name`.\n//\n// The text lives inside a single run that reads (decoded):\n//   ...caseManagementLocation.venue_name>><<else>> Online Civil Claims<<es_>>\n// We only need to swap the `venue_name` token for `external_short_name`,\n// keeping the surrounding `<<...>>` template markup and formatting intact.\n\nconst results = context.document.body.search(\"venue_name\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"external_short_name\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The SDO template merge-field referencing the case management location's\n# venue name field should instead reference the location ref data's\n# \"external short name\" field: `venue_name` -> `external_short_name`.\n#\n# The text lives inside a single run that reads (decoded):\n#   ...caseManagementLocation.venue_name>><<else>> Online Civil Claims<<es_>>\n# We only need to swap the `venue_name` token for `external_short_name`,\n# keeping the surrounding `<<...>>` template markup and formatting intact.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"venue_name\"\n$find.Replacement.Text = \"external_short_name\"\n$find.Execute(\"venue_name\", $false, $false, $false, $false, $false, $true, 1, $false, \"external_short_name\", 2)\n"}
